# Wire Harness / ConnectorProfile.xlsx update
# "Updated Connector for Front/Rear Brake Pressure and 6DOF IMU"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Front Brake Pressure (row 6, CN-005) ---
# Fill in the newly-chosen connector part / manufacturer / position count.
$ws.Range("C6").Value = "ASL 0-06-05SC-HE"
$ws.Range("D6").Value = "Deutsch"
$ws.Range("E6").Value = 3

# --- Rear Brake Pressure (row 7, CN-006) ---
$ws.Range("C7").Value = "ASL 0-06-05SC-HE"
$ws.Range("D7").Value = "Deutsch"
$ws.Range("E7").Value = 3

# --- Front Accelerometer (row 8, CN-007) : number of positions filled in ---
$ws.Range("E8").Value = 5

# --- Rear Accelerometer (row 9, CN-008) : number of positions filled in ---
$ws.Range("E9").Value = 5

# --- IMU (row 15, CN-014) : connector part number updated for 6DOF IMU ---
$ws.Range("C15").Value = "ASL606-05SN"

# Restore the view: selection on C6, scrolled back to the top of the sheet.
$ws.Range("A1").Select() | Out-Null
$ws.Range("C6").Select() | Out-Null
